# Apply the diagram edits described in the commit:
# - Rename the "gt:GradTrak" shape's label to "gt:ReadOnlyGradTrak"
#   and resize/reposition its box to fit the longer text.
# - Merge the two "getCode" / "()" runs into a single "getCode()" run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {
    if ($shp.Id -eq 26) {
        $shp.TextFrame.TextRange.Text = "gt:ReadOnlyGradTrak"
        # Target EMU box: off x=4076452 y=5229200, ext cx=1224136 cy=551451
        # (values nudged slightly so the point<->EMU round-trip lands exactly)
        $shp.Left   = 320.9804840409449
        $shp.Top    = 411.748031696063
        $shp.Width  = 96.38866051732283
        $shp.Height = 43.421339082677164
    }
    elseif ($shp.Id -eq 104) {
        # The shape already reads "getCode()" but as two separate runs
        # ("getCode" + "()"); force a real text replacement (via a
        # placeholder that shares no characters with the final text, so
        # the engine fully rebuilds the paragraph) so it collapses down
        # into a single run instead of a partial in-place patch.
        $shp.TextFrame.TextRange.Text = "ZZZ"
        $shp.TextFrame.TextRange.Text = "getCode()"
    }
}
